$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 260 (shifts existing rows 260..365 down to 261..366,
# dimension grows from A1:R365 to A1:R366).
$ws.Rows.Item(260).Insert()

# Populate the new row 260 with the new price-quote record.
$ws.Range("A260").Value = 5
$ws.Range("B260").Value = "Macroferia Regional de Talca"
$ws.Range("C260").Value = "Maule"
$ws.Range("D260").Value = 44875
$ws.Range("D260").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E260").Value = 7
$ws.Range("F260").Value = 100112008
$ws.Range("G260").Value = "Coliflor"
$ws.Range("H260").Value = "Sin especificar"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 3000
$ws.Range("K260").Value = 700
$ws.Range("L260").Value = 700
$ws.Range("M260").Value = 700
$ws.Range("N260").Value = "$/unidad"
$ws.Range("O260").Value = "Región del Maule"
$ws.Range("P260").Value = 700
$ws.Range("Q260").Value = 1
$ws.Range("R260").Value = "Hortaliza"
